# Planning View Generator - initial build
# ------------------------------------------------------------------
# 1) items: drop the safetyStock column (old column C) - strategy,
#    currentInventory, FOQ and palletQty all shift one column left.
# 2) resources: add a new "class" column classifying each resource as
#    Machine or Labor.
# 3) constraints: drop the computed "days" helper column (old column D,
#    formula =C/20) - no longer needed.
# 4) Selections / active sheet are updated to match the new layout.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- items sheet: remove safetyStock column -----------------------
$wsItems = $wb.Worksheets.Item("items")
$wsItems.Activate()
$wsItems.Columns("C").Delete()

# --- resources sheet: add class column -----------------------------
$wsRes = $wb.Worksheets.Item("resources")
$wsRes.Activate()

$wsRes.Range("C1").Value = "class"
$wsRes.Range("A1").Copy()
$wsRes.Range("C1").PasteSpecial(-4122)

$wsRes.Range("C2").Value = "Machine"
$wsRes.Range("C3").Value = "Labor"
$wsRes.Range("C4").Value = "Labor"
$wsRes.Range("C5").Value = "Machine"

# --- constraints sheet: remove the computed days column -------------
$wsCon = $wb.Worksheets.Item("constraints")
$wsCon.Activate()
$wsCon.Columns("D").Delete()
$wsCon.Range("D2:D13").Select()
$wsCon.Range("E2").Select()

# --- calendar sheet: selection only ---------------------------------
$wsCal = $wb.Worksheets.Item("calendar")
$wsCal.Activate()
$wsCal.Range("H12").Select()

# --- restore per-sheet selections ------------------------------------
$wsItems.Activate()
$wsItems.Range("L3").Select()

$wsRes.Activate()
$wsRes.Range("D3").Select()

# constraints ends up the active tab, matching the saved workbook view
$wsCon.Activate()
$wsCon.Range("E2").Select()
